$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.152819514274597
$ws.Range("B1").Value = 1.432760715484619
$ws.Range("C1").Value = 6.712436676025391
$ws.Range("D1").Value = 2.11025857925415
$ws.Range("E1").Value = 0.9321473836898804
